# The commit adds three blank paragraphs followed by a new paragraph of
# text ("Juan de Dios Pantajas pruebas de git ") at the end of the
# document body, right after the existing single paragraph.
#
# The new blank paragraphs keep only the inherited `lang="en-US"` mark
# (no runs at all), and the final paragraph's runs carry no rPr and are
# interrupted by spell-check <w:proofErr> markers around "Pantajas" and
# "git" - exactly mirroring a live Word typing+proofing session rather
# than a single formatted insert. We reproduce that exact OOXML shape
# with Range.InsertXML so no extra/implicit formatting gets injected by
# a higher-level Selection/Range.Text insert.

$d = $word.ActiveDocument

# Collapsed range at the very end of the document's main story.
$endOfDoc = $d.Content.End - 1
$r = $d.Range($endOfDoc, $endOfDoc)

$bodyFragment = (
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' +
  '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' +
  '<w:p>' +
    '<w:r><w:t xml:space="preserve">Juan de Dios </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Pantajas</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> pruebas d</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">e </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>git</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>'
)

$packageXml = (
  '<?xml version="1.0" standalone="yes"?>' +
  '<?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" ' +
      'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $bodyFragment + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
)

$r.InsertXML($packageXml)
